$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = 0.2192982456140351
$ws.Range("C2").Value = 0.5131578947368421
$ws.Range("J2").Value = 0.02631578947368421
$ws.Range("O2").Value = 0.004385964912280702
$ws.Range("P2").Value = 0.1491228070175439
$ws.Range("S2").Value = 0.08771929824561403
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.05555555555555555
$ws.Range("J3").Value = 0.007936507936507936
$ws.Range("P3").Value = 0.6507936507936508
$ws.Range("S3").Value = 0.2777777777777778
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.5853658536585366
$ws.Range("S4").Value = 0.3414634146341464
$ws.Range("B6").Value = 0.06
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.055
$ws.Range("J6").Value = 0.145
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.09
$ws.Range("S6").Value = 0.505
$ws.Range("B7").Value = 0.07291666666666667
$ws.Range("D7").Value = 0.005208333333333333
$ws.Range("F7").Value = 0.02604166666666667
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.02604166666666667
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.08854166666666667
$ws.Range("S7").Value = 0.4479166666666667
$ws.Range("B8").Value = 0.05510534846029173
$ws.Range("D8").Value = 0.01782820097244733
$ws.Range("E8").Value = 0.001620745542949757
$ws.Range("F8").Value = 0.04862236628849271
$ws.Range("J8").Value = 0.1442463533225284
$ws.Range("O8").Value = 0.009724473257698542
$ws.Range("Q8").Value = 0.1977309562398703
$ws.Range("R8").Value = 0.06969205834683954
$ws.Range("S8").Value = 0.4554294975688817
$ws.Range("B9").Value = 0.0663265306122449
$ws.Range("D9").Value = 0.00510204081632653
$ws.Range("F9").Value = 0.07653061224489796
$ws.Range("J9").Value = 0.1224489795918367
$ws.Range("O9").Value = 0.01020408163265306
$ws.Range("Q9").Value = 0.1683673469387755
$ws.Range("R9").Value = 0.05612244897959184
$ws.Range("S9").Value = 0.4948979591836735
$ws.Range("B10").Value = 0.07236842105263158
$ws.Range("D10").Value = 0.01608187134502924
$ws.Range("E10").Value = 0.002923976608187134
$ws.Range("F10").Value = 0.07163742690058479
$ws.Range("J10").Value = 0.1264619883040936
$ws.Range("O10").Value = 0.01608187134502924
$ws.Range("Q10").Value = 0.2134502923976608
$ws.Range("R10").Value = 0.07456140350877193
$ws.Range("S10").Value = 0.4064327485380117
$ws.Range("G11").Value = 0.1117824773413897
$ws.Range("J11").Value = 0.1148036253776435
$ws.Range("K11").Value = 0.1873111782477341
$ws.Range("L11").Value = 0.5619335347432024
$ws.Range("S11").Value = 0.02416918429003021
$ws.Range("G12").Value = 0.6736842105263158
$ws.Range("J12").Value = 0.2947368421052631
$ws.Range("L12").Value = 0.01578947368421053
$ws.Range("S12").Value = 0.01578947368421053
$ws.Range("G13").Value = 0.7727272727272727
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("F15").Value = 0.004716981132075472
$ws.Range("H15").Value = 0.1792452830188679
$ws.Range("I15").Value = 0.06132075471698113
$ws.Range("J15").Value = 0.3726415094339622
$ws.Range("K15").Value = 0.08490566037735849
$ws.Range("M15").Value = 0.01415094339622642
$ws.Range("O15").Value = 0.04716981132075472
$ws.Range("S15").Value = 0.2358490566037736
$ws.Range("F16").Value = 0.007246376811594203
$ws.Range("H16").Value = 0.2391304347826087
$ws.Range("I16").Value = 0.05797101449275362
$ws.Range("J16").Value = 0.3695652173913043
$ws.Range("K16").Value = 0.1449275362318841
$ws.Range("M16").Value = 0.007246376811594203
$ws.Range("N16").Value = 0.007246376811594203
$ws.Range("O16").Value = 0.02173913043478261
$ws.Range("S16").Value = 0.1449275362318841
$ws.Range("F17").Value = 0.01183431952662722
$ws.Range("H17").Value = 0.2169625246548323
$ws.Range("I17").Value = 0.07495069033530571
$ws.Range("J17").Value = 0.4102564102564102
$ws.Range("K17").Value = 0.09467455621301775
$ws.Range("M17").Value = 0.01577909270216963
$ws.Range("O17").Value = 0.05128205128205128
$ws.Range("S17").Value = 0.1242603550295858
$ws.Range("F18").Value = 0.005235602094240838
$ws.Range("H18").Value = 0.2041884816753927
$ws.Range("I18").Value = 0.06282722513089005
$ws.Range("J18").Value = 0.4240837696335079
$ws.Range("K18").Value = 0.1047120418848168
$ws.Range("M18").Value = 0.01570680628272251
$ws.Range("O18").Value = 0.05759162303664921
$ws.Range("S18").Value = 0.1256544502617801
$ws.Range("F19").Value = 0.007822685788787484
$ws.Range("H19").Value = 0.2620599739243807
$ws.Range("I19").Value = 0.08148631029986962
$ws.Range("J19").Value = 0.3376792698826597
$ws.Range("K19").Value = 0.1010430247718383
$ws.Range("M19").Value = 0.01890482398956975
$ws.Range("N19").Value = 0.000651890482398957
$ws.Range("O19").Value = 0.06127770534550196
$ws.Range("S19").Value = 0.1290743155149935
